$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52
$ws.Range("A52").Value = 112105381
$ws.Range("B52").Value = 88899
$ws.Range("E52").Value = 3286
$ws.Range("F52").Value = "Flattoppad klubbsvamp"
$ws.Range("G52").Value = "Clavariadelphus truncatus"
$ws.Range("H52").Value = "(Quél.) Donk"
$ws.Range("P52").Value = "Landverktjärnen (Landverktjärnen), Jmt"
$ws.Range("Q52").Value = 446564
$ws.Range("R52").Value = 7032716
$ws.Range("Z52").ClearContents()
$ws.Range("AB52").ClearContents()

# Row 53
$ws.Range("A53").Value = 112110532
$ws.Range("B53").Value = 88002
$ws.Range("D53").Value = "VU"
$ws.Range("E53").Value = 245031
$ws.Range("F53").Value = "Borgsjömusseron"
$ws.Range("G53").Value = "Tricholoma borgsjoeënse"
$ws.Range("H53").Value = "Jacobsson & Muskos"
$ws.Range("P53").Value = "Svensbergsbäcken (Svensbergsbäcken), Jmt"
$ws.Range("Q53").Value = 446765
$ws.Range("R53").Value = 7032863
$ws.Range("AC53").Value = "Längs med en liten stig. På svag sluttning bland kam-och husmossa, revlummer och ekbräken. I närheten finns granvaxskivling, rosa/besk vaxskivling och äggvaxskivling."
$ws.Range("Z53").ClearContents()
$ws.Range("AB53").ClearContents()

# Row 54
$ws.Range("A54").Value = 112111388
$ws.Range("B54").Value = 88956
$ws.Range("D54").Value = "VU"
$ws.Range("E54").Value = 5747
$ws.Range("F54").Value = "Läderdoftande fingersvamp"
$ws.Range("G54").Value = "Ramaria safraniolens"
$ws.Range("H54").Value = "Christian"
$ws.Range("Q54").Value = 446734
$ws.Range("R54").Value = 7032709
$ws.Range("Z54").ClearContents()
$ws.Range("AB54").ClearContents()

# Row 55
$ws.Range("A55").Value = 112104863
$ws.Range("B55").Value = 90651
$ws.Range("D55").Value = "NT"
$ws.Range("E55").Value = 1968
$ws.Range("F55").Value = "Grantaggsvamp"
$ws.Range("G55").Value = "Bankera violascens"
$ws.Range("H55").Value = "(Alb. & Schwein. : Fr.) Pouzar"
$ws.Range("P55").Value = "Hökån (Hökån), Jmt"
$ws.Range("Q55").Value = 446637
$ws.Range("R55").Value = 7032524
$ws.Range("Z55").ClearContents()
$ws.Range("AB55").ClearContents()

# Row 56
$ws.Range("A56").Value = 112111378
$ws.Range("B56").Value = 82949
$ws.Range("D56").Value = "NT"
$ws.Range("E56").Value = 5589
$ws.Range("F56").Value = "Rödbrun klubbdyna"
$ws.Range("G56").Value = "Trichoderma nybergianum"
$ws.Range("H56").Value = "(T.Ulvinen & H.L.Chamb.) Jaklitsch & Voglmayr"
$ws.Range("P56").Value = "Renkullmyren (Renkullmyren), Jmt"
$ws.Range("Q56").Value = 446760
$ws.Range("R56").Value = 7032715
$ws.Range("Z56").ClearContents()
$ws.Range("AB56").ClearContents()
$ws.Range("AC56").ClearContents()

# Row 57
$ws.Range("A57").Value = 112111398
$ws.Range("B57").Value = 88966
$ws.Range("E57").Value = 5754
$ws.Range("F57").Value = "Gultoppig fingersvamp"
$ws.Range("G57").Value = "Ramaria testaceoflava"
$ws.Range("H57").Value = "(Bres.) Corner"
$ws.Range("P57").Value = "Renkullmyren (Renkullmyren), Jmt"
$ws.Range("Q57").Value = 446740
$ws.Range("R57").Value = 7032705
$ws.Range("Z57").ClearContents()
$ws.Range("AB57").ClearContents()

# Row 58
$ws.Range("A58").Value = 112111386
$ws.Range("B58").Value = 88946
$ws.Range("E58").Value = 256335
$ws.Range("F58").Value = "Taggfingersvamp"
$ws.Range("G58").Value = "Ramaria karstenii"
$ws.Range("H58").Value = "(Sacc. & P.Syd.) Corner"
$ws.Range("Q58").Value = 446734
$ws.Range("R58").Value = 7032709
$ws.Range("Z58").ClearContents()
$ws.Range("AB58").ClearContents()

# Row 59
$ws.Range("A59").Value = 112213232
$ws.Range("B59").Value = 88956
$ws.Range("D59").Value = "VU"
$ws.Range("E59").Value = 5747
$ws.Range("F59").Value = "Läderdoftande fingersvamp"
$ws.Range("G59").Value = "Ramaria safraniolens"
$ws.Range("H59").Value = "Christian"
$ws.Range("Q59").Value = 446675
$ws.Range("R59").Value = 7032593

# Row 60
$ws.Range("A60").Value = 112213255
$ws.Range("B60").Value = 90332
$ws.Range("D60").Value = "LC"
$ws.Range("E60").Value = 4769
$ws.Range("F60").Value = "Svavelriska"
$ws.Range("G60").Value = "Lactarius scrobiculatus"
$ws.Range("H60").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q60").Value = 446605
$ws.Range("R60").Value = 7032710
